$d = $word.ActiveDocument

# 1. Title "Contabilidad Casa": bump font size from 18pt (sz/szCs=36) to 24pt (sz/szCs=48).
#    Setting Font.Size / Font.SizeBi on the paragraph's Range updates both the run's
#    run-properties and the paragraph-mark run-properties (w:pPr/w:rPr), matching
#    the diff which changes both occurrences.
$title = $d.Paragraphs(1).Range
$title.Font.Size = 24
$title.Font.SizeBi = 24

# 2. Move the "_GoBack" bookmark so it sits right after the "Contabilidad Casa" run
#    (it used to sit right after "txbMontofD" further down in the document).
#    A zero-length Range placed exactly at a paragraph's text end (one position
#    before the paragraph mark) can't be targeted directly with Bookmarks.Add, so
#    temporarily append a marker character, anchor the bookmark right before the
#    marker, then delete the marker again - leaving the bookmark cleanly after the
#    run, still inside the paragraph.
$title = $d.Paragraphs(1).Range
$insertPos = $title.End
$title.InsertAfter("#")

$markerStart = $insertPos - 1
$markerEnd = $insertPos
$anchor = $d.Range($markerStart, $markerStart)

# Adding a bookmark named "_GoBack" automatically relocates it - Word only ever
# keeps a single "_GoBack" bookmark - so this also removes the stale one that used
# to sit after "txbMontofD".
$d.Bookmarks.Add("_GoBack", $anchor)

$marker = $d.Range($markerStart, $markerEnd)
$marker.Delete()
